# https://jira.hl7.org/browse/FHIR-49085
# In diagram, corrected "Payer CDA System" to "Payer CDS System"

$oldText = "Payer CDA / Production Systems"
$newText = "Payer CDS / Production Systems"

function Update-ShapeText($shp) {
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq $oldText) {
                $shp.TextFrame.TextRange.Text = $newText
            }
        }
    }
    if ($shp.Type -eq 6) {
        # msoGroup: walk every item inside the group too.
        $items = $shp.GroupItems
        for ($i = 1; $i -le $items.Count; $i++) {
            Update-ShapeText $items.Item($i)
        }
    }
}

$p = $ppt.ActivePresentation
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        Update-ShapeText $s.Shapes.Item($shi)
    }
}
